# "Changed all links to Ryzen 7000"
# Rename the tracked CPUs to the Ryzen 7000 lineup, refresh their prices
# and check-date, and add a new row for the Ryzen 9 7950X.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkDate = 44999.415338126288

# Row 2 -> AMD Ryzen 5 7600X
$ws.Range("A2").Value = "AMD Ryzen 5 7600X"
$ws.Range("B2").Value = 249.99
$ws.Range("C2").Value = 319.99
$ws.Range("D2").Value = $checkDate

# Row 3 -> AMD Ryzen 7 7700X
$ws.Range("A3").Value = "AMD Ryzen 7 7700X"
$ws.Range("B3").Value = 344.99
$ws.Range("C3").Value = 444.99
$ws.Range("D3").Value = $checkDate

# Row 4 -> AMD Ryzen 9 7900X
$ws.Range("A4").Value = "AMD Ryzen 9 7900X"
$ws.Range("B4").Value = 436.55
$ws.Range("C4").Value = 579.99
$ws.Range("D4").Value = $checkDate

# Row 5 (new) -> AMD Ryzen 9 7950X
$ws.Range("A5").Value = "AMD Ryzen 9 7950X"
$ws.Range("B5").Value = 599.99
$ws.Range("C5").Value = 759.99
$ws.Range("D5").Value = $checkDate

# New D5 needs the same date number format as the rest of column D.
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$ws.Range("D5").Value = $checkDate

# Column A auto-shrinks now that the CPU names are shorter.
$ws.Range("A1").EntireColumn.ColumnWidth = 16.5
